# Apply updated boolean vote values to the threshold approval matrix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = $true;  "C2" = $true;  "D2" = $true;  "E2" = $true;  "F2" = $true;  "G2" = $true;
    "D3" = $true;  "F3" = $true;
    "G4" = $false;
    "B5" = $true;  "C5" = $false; "D5" = $true;  "E5" = $true;  "G5" = $false;
    "C6" = $true;  "D6" = $true;  "E6" = $true;  "F6" = $true;  "G6" = $true;
    "D7" = $true;  "F7" = $false;
    "B10" = $true;
    "F11" = $false;
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
